# đôi font-end mới và thêm pdf
# Insert a new leading "STT" (sequence number) column into the awards table,
# shifting the existing "Đơn Vị" / "Năm" / award columns one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new column at position A; this shifts all existing
#    columns (and their data/styles/widths) one slot to the right.
$ws.Columns.Item(1).Insert()

# 2) Give the new column header cell (A1) the same look as the other
#    header cells (bold / fill / border) by copying the format from B1,
#    then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "STT"

# 3) Give the new data cells (A2:A10) the same look as the other data
#    cells by copying the format from B2, then fill in the sequence
#    numbers 1..9.
$ws.Range("B2").Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9

$excel.CutCopyMode = 0

# 4) Narrow the new STT column to fit its short content (close to the
#    6.856-character "best fit" width used by the rest of the sheet).
$ws.Columns.Item(1).ColumnWidth = 6

# 5) Refresh the sheet's remembered selection so it spans the full,
#    now-wider table (A2:J10) instead of the old A2:I10 range.
$ws.Range("A2:J10").Select() | Out-Null
